$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "RCVQ" -> "log RCVQ" (TextBox 15) — right-aligned autofit textbox grows leftward
$shapeQ = $s.Shapes.Item(11)
$shapeQ.TextFrame.TextRange.Text = "log RCVQ"
$shapeQ.Left = 266.9637145996094
$shapeQ.Width = 61.62102508544922

# "RCVM" -> "log RCVM" (TextBox 16) — right-aligned autofit textbox grows leftward
$shapeM = $s.Shapes.Item(12)
$shapeM.TextFrame.TextRange.Text = "log RCVM"
$shapeM.Left = 266.33270263671875
$shapeM.Width = 62.25204849243164
